# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-04-03 Thursday", 2)

# Update each table cell by exact (row, column) position so the two
# duplicate "482÷2=241, 0" cells get their correct, distinct replacements.
$t = $word.ActiveDocument.Tables.Item(1)

$cells = @(
    @{Row=1;  Col=1; New="879÷6=146, 3"},
    @{Row=1;  Col=2; New="726÷9=80, 6"},
    @{Row=1;  Col=3; New="770÷5=154, 0"},
    @{Row=1;  Col=4; New="417÷8=52, 1"},
    @{Row=1;  Col=5; New="414÷5=82, 4"},

    @{Row=5;  Col=1; New="640÷5=128, 0"},
    @{Row=5;  Col=2; New="580÷2=290, 0"},
    @{Row=5;  Col=3; New="105÷3=35, 0"},
    @{Row=5;  Col=4; New="520÷4=130, 0"},
    @{Row=5;  Col=5; New="834÷9=92, 6"},

    @{Row=9;  Col=1; New="956÷4=239, 0"},
    @{Row=9;  Col=2; New="558÷3=186, 0"},
    @{Row=9;  Col=3; New="542÷3=180, 2"},
    @{Row=9;  Col=4; New="836÷3=278, 2"},
    @{Row=9;  Col=5; New="763÷4=190, 3"},

    @{Row=13; Col=1; New="522÷9=58, 0"},
    @{Row=13; Col=2; New="763÷9=84, 7"},
    @{Row=13; Col=3; New="166÷9=18, 4"},
    @{Row=13; Col=4; New="696÷3=232, 0"},
    @{Row=13; Col=5; New="614÷2=307, 0"},

    @{Row=17; Col=1; New="787÷8=98, 3"},
    @{Row=17; Col=2; New="434÷6=72, 2"},
    @{Row=17; Col=3; New="801÷3=267, 0"},
    @{Row=17; Col=4; New="740÷6=123, 2"},
    @{Row=17; Col=5; New="412÷6=68, 4"}
)

foreach ($item in $cells) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, preserving the cell's paragraph formatting.
    $rng.End = $rng.End - 1
    $rng.Text = $item.New
}
